# Update mailing address in terms of use (and "Last Updated" date)
# https://phabricator.endlessm.com/T17092
#
# Old address: 512 2nd Street, Floor 3 / Third Floor, San Francisco, CA 94107
# New address: 575 Market Street, Suite 825, San Francisco, CA 94105

$d = $word.ActiveDocument
$tab = [char]9

# 1) Standalone address line (non-bold), e.g. "...512 2nd Street, Floor 3"
$ok1 = $d.Content.Find.Execute(
    "512 2nd Street, Floor 3", $true, $false, $false, $false, $false,
    $true, 1, $false, "575 Market Street, Suite 825", 2)

# 2) The city/zip line immediately following it: "San Francisco, CA 94107" -> "...94105"
#    (scoped with the preceding tabs so we only touch this standalone line)
$ok2 = $d.Content.Find.Execute(
    "$tab$tab" + "San Francisco, CA 94107", $true, $false, $false, $false, $false,
    $true, 1, $false, "$tab$tab" + "San Francisco, CA 94105", 2)

# 3) Inline (regular weight) mention: "...512 2nd Street, Third Floor, San Francisco, CA 94107."
$ok3 = $d.Content.Find.Execute(
    "512 2nd Street, Third Floor, San Francisco, CA 94107.", $true, $false, $false, $false, $false,
    $true, 1, $false, "575 Market Street, Suite 825, San Francisco, CA 94105.", 2)

# 4) Inline (bold) mention: "...512 2nd Street, Third Floor, San Francisco, CA 94107 " (trailing
#    space keeps the match inside the bold run, not spilling into the following "dengan..." run)
$ok4 = $d.Content.Find.Execute(
    "512 2nd Street, Third Floor, San Francisco, CA 94107 ", $true, $false, $false, $false, $false,
    $true, 1, $false, "575 Market Street, Suite 825, San Francisco, CA 94105 ", 2)

# 5) "Last Updated" date in the page header: 15 Februari 2017 -> 2 Juni 2017
$sec = $d.Sections(1)
$hdr = $sec.Headers(1)
$ok5 = $hdr.Range.Find.Execute(
    "Pembaruan Terakhir: 15 Februari 2017", $true, $false, $false, $false, $false,
    $true, 1, $false, "Pembaruan Terakhir: 2 Juni 2017", 2)

Write-Output "address line 1 replaced: $ok1"
Write-Output "zip line replaced: $ok2"
Write-Output "inline regular mention replaced: $ok3"
Write-Output "inline bold mention replaced: $ok4"
Write-Output "header date replaced: $ok5"
